# Updates the cryptos price list (column D = Price, column E = Volume(1h)).
# Also fixes the 3-way reorder of Quant / TrustWalletToken / Aave in rows 42-44,
# which changes every cell (Coin name, Link, Price, Volume) in those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking price/percent strings to stay TEXT (matching the
# original inline-string cells) instead of Excel auto-converting them to numbers.
$dataRange = $ws.Range('D2:E51')
$dataRange.NumberFormat = '@'

$ws.Range('D2').Value = '29.390.92'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '1.876.71'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('D5').Value = '0.7173'
$ws.Range('E5').Value = '  +1.08%  '
$ws.Range('D6').Value = '243.72'
$ws.Range('E6').Value = '  +0.70%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.07941'
$ws.Range('E8').Value = '  +1.38%  '
$ws.Range('D9').Value = '0.3145'
$ws.Range('E9').Value = '  +1.07%  '
$ws.Range('D10').Value = '24.96'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').Value = '0.08130'
$ws.Range('E11').Value = '  -3.19%  '
$ws.Range('D12').Value = '1.886.21'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '95.69'
$ws.Range('E13').Value = '  +4.86%  '
$ws.Range('D14').Value = '5.236'
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').Value = '0.7067'
$ws.Range('E15').Value = '  -1.64%  '
$ws.Range('D16').Value = '6.397'
$ws.Range('E16').Value = '  +4.20%  '
$ws.Range('D17').Value = '0.000008426'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').Value = '29.397.03'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').Value = '252.89'
$ws.Range('E19').Value = '  +5.05%  '
$ws.Range('D20').Value = '13.40'
$ws.Range('E20').Value = '  +1.40%  '
$ws.Range('D21').Value = '2.138.57'
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').Value = '7.671'
$ws.Range('E23').Value = '  -1.03%  '
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').Value = '0.1589'
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').Value = '9.072'
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('D27').Value = '161.96'
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('D28').Value = '18.91'
$ws.Range('E28').Value = '  +2.05%  '
$ws.Range('D29').Value = '1.506'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '4.415'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').Value = '4.297'
$ws.Range('E31').Value = '  -1.07%  '
$ws.Range('D32').Value = '1.222'
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('D33').Value = '0.05324'
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').Value = '0.7567'
$ws.Range('E35').Value = '  +1.28%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('D38').Value = '0.01892'
$ws.Range('E38').Value = '  +0.67%  '
$ws.Range('D39').Value = '1.265.87'
$ws.Range('E39').Value = '  +1.78%  '
$ws.Range('E40').Value = '  +0.92%  '
$ws.Range('D41').Value = '6.390'
$ws.Range('E41').Value = '  -1.84%  '
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('D47').Value = '2.035.67'
$ws.Range('E47').Value = '  +0.71%  '
$ws.Range('D48').Value = '1.810'
$ws.Range('E48').Value = '  +0.55%  '
$ws.Range('D49').Value = '0.5202'
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('D50').Value = '9.518'
$ws.Range('E50').Value = '  +0.80%  '
$ws.Range('D51').Value = '0.4347'
$ws.Range('E51').Value = '  -0.01%  '

# Rows 42-44: Quant / TrustWalletToken / Aave rotate order with new data.
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.9060'
$ws.Range('E42').Value = '  +1.51%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '74.33'
$ws.Range('E43').Value = '  +2.59%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '111.61'
$ws.Range('E44').Value = '  +1.23%  '

# Restore default (General/Normal) styling so no stray number-format/style
# artifacts are left behind on the cells we touched.
$dataRange.Style = 'Normal'
